# Update Pinlayout spreadsheet with keypad pin mapping
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sender")

# Map keypad pin labels onto the FUNCTION column (G) for PORT E pins 2,3,6,7,8,9,10,11
$ws.Range("G21").Value = "Keypad Pin 5"
$ws.Range("G22").Value = "Keypad Pin 6"
$ws.Range("G25").Value = "Keypad Pin 7"
$ws.Range("G26").Value = "Keypad Pin 8"
$ws.Range("G27").Value = "Keypad Pin 1"
$ws.Range("G28").Value = "Keypad Pin 2"
$ws.Range("G29").Value = "Keypad Pin 3"
$ws.Range("G30").Value = "Keypad Pin 4"

# Update the frozen-pane view / active selection to where the user left off
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$ws.Range("J26").Select()
